$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-cell value updates (price/volume refresh); rows 40/41 content swapped

$ws.Range("D2").Value = "26.025.94"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "1.644.14"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("E4").Value = "  +0.56%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.01"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("E7").Value = "  +0.50%  "
$ws.Range("E8").Value = "  +0.45%  "
$ws.Range("E9").Value = "  +0.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.56"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0798"
$ws.Range("D11").ClearFormats()
$ws.Range("E12").Value = "  +0.68%  "
$ws.Range("D13").Value = "1.650.52"
$ws.Range("E13").Value = "  +0.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.544"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.44"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.39%  "
$ws.Range("D17").Value = "26.058.22"
$ws.Range("E17").Value = "  +0.41%  "
$ws.Range("E18").Value = "  +0.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "194.30"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.28%  "
$ws.Range("E20").Value = "  -0.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.93"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("E22").Value = "  -1.05%  "
$ws.Range("E23").Value = "  +3.84%  "
$ws.Range("E24").Value = "  -0.70%  "
$ws.Range("E25").Value = "  +0.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "143.24"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.43%  "
$ws.Range("E27").Value = "  +0.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.51"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.47%  "
$ws.Range("E29").Value = "  +0.45%  "
$ws.Range("E30").Value = "  -1.08%  "
$ws.Range("E31").Value = "  -0.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.27"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.43%  "
$ws.Range("E33").Value = "  -0.59%  "
$ws.Range("E34").Value = "  +1.62%  "
$ws.Range("E36").Value = "  -0.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.540"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.94%  "
$ws.Range("E38").Value = "  +0.33%  "
$ws.Range("E39").Value = "  -0.12%  "
$ws.Range("B40").Value = "Quant"
$ws.Range("C40").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "99.18"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.19%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.46"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.74%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("E43").Value = "  +2.33%  "
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("E45").Value = "  -1.44%  "
$ws.Range("E46").Value = "  +2.64%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.77"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.42%  "
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("E49").Value = "  +0.42%  "
$ws.Range("E50").Value = "  -0.94%  "
$ws.Range("E51").Value = "  +3.12%  "
